$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 and 8: USDC/XRP swap ranking positions
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"

# Column D (Price) updates
$ws.Range("D2").Value = "65.640.62"
$ws.Range("D3").Value = "2.660.56"
$ws.Range("D5").Value = "'598.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'159.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.640"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.126"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'5.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.398"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = "'29.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000194"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.143.71"
$ws.Range("D16").Value = "65.609.85"
$ws.Range("D17").Value = "2.659.77"
$ws.Range("D18").Value = "'12.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'4.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'353.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'7.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'69.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'1.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.0000113"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'1.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'566.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").Value = "'2.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'1.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'6.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'5.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'20.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'1.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'154.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'2.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'161.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'4.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.0616"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'23.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.644"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.0258"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.101"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'19.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "0.0₆0244"
$ws.Range("D51").Value = "'0.814"
$ws.Range("D51").Style = "Normal"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +3.86%  "
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +9.20%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("E28").Value = "  +6.67%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  -3.27%  "
$ws.Range("E41").Value = "  +6.49%  "
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("E50").Value = "  -7.08%  "
$ws.Range("E51").Value = "  +0.45%  "
